$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "Rectangle 3": "Not endowed" -> "Less endowed" ---
$sh3 = $s.Shapes.Item("Rectangle 3")
$tr3 = $sh3.TextFrame.TextRange

# Drop the leading "Not" (3 characters). What's left (" endowed") starts
# with the space run, which was never explicitly bold, so rebuilding the
# wording through it keeps the final single run non-bold without having to
# call Font.Bold explicitly (which would otherwise stamp an explicit b="0").
$tr3.Characters(1, 3).Text = ""
# Retype through a throwaway one-word placeholder first: replacing like-for-
# like word counts preserves the old run split, but going through a single
# "odd" word collapses the paragraph down to one run. Then set real text.
$tr3.Text = "X"
$tr3.Text = "Less endowed"
$tr3.Font.Size = 16

# Resize/reposition after editing the text -- this is an autofit text box,
# so PowerPoint recomputes its height from the text on every edit and any
# earlier size/position change would just get overwritten.
$sh3.Top = 311.9234
$sh3.Height = 26.6578

# --- "Rectangle 4": "Endowed" -> "More endowed" ---
$sh4 = $s.Shapes.Item("Rectangle 4")
$tr4 = $sh4.TextFrame.TextRange

$tr4.Text = "More endowed"
$tr4.Font.Bold = $false
$tr4.Font.Size = 16

$sh4.Top = 163.89921
$sh4.Height = 26.6578
